# Apply "progress as of date 04 Nov 2025" update to the Training Dashboard sheet.
# For each data row (3-17): decrement the "PERIOD TO EXPIRE" (col H) by 1,
# and bump the "LAST UPDATE" (col I) text from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 17; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H: PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    # Column I: LAST UPDATE - keep it as literal text (not auto-converted to a
    # date serial) by forcing a text number format before assigning the value.
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value2 = "04-Nov-2025"
}
